$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Complete the existing "Day 9" row (row 10) ---
$ws.Range("C10").Value = "Longest Repeating Character Replacement"
$ws.Range("D10").Value = "Majority Element"
$ws.Range("E10").Value = "Rotate Array"
$ws.Range("F10").Value = "Sliding Window, Arrays, HashMap"
$ws.Range("G10").Value = "S"
$ws.Range("H10").Value = "YES"

# --- Day labels for the two new rows ---
$ws.Range("A11").Value = "Day 10"
$ws.Range("A12").Value = "Day 11"

# --- Day 10 (row 11) ---
$ws.Range("B11").NumberFormat = "m/d/yy"
$ws.Range("B11").Value = 45812
$ws.Range("C11").Value = "Happy Number"
$ws.Range("D11").Value = "Ransom Note"
$ws.Range("E11").Value = "Word Pattern"
$ws.Range("F11").Value = "HashSet, HashMap, String Matching"
$ws.Range("G11").Value = "S"
$ws.Range("H11").Value = "YES"

# --- Day 11 (row 12) ---
$ws.Range("B12").NumberFormat = "m/d/yy"
$ws.Range("B12").Value = 45813
$ws.Range("C12").Value = "Is Subsequence"
$ws.Range("D12").Value = "Roman to Integer"
$ws.Range("E12").Value = "Integer to Roman"
$ws.Range("F12").Value = "String Matching, Two Pointers, Math"
$ws.Range("G12").Value = "S"
$ws.Range("H12").Value = "YES"

# --- Column C widened to fit the longer text ---
$ws.Columns.Item(3).ColumnWidth = 32

# --- Selection moves to A13 (the next empty row) after the edits ---
$ws.Range("A13").Select()
